$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the placeholder " " first_shift value is cleared (no shift detected within 500ms gap-corrected window)
$clearRows = @(2,3,4,6,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,45,46,48,54,56,61,62,63,64,70,71,73,74,76,77,78)
foreach ($r in $clearRows) {
    $ws.Cells.Item($r, 6).ClearContents()
}

# Rows with recomputed shift classification / latency / counts after allowing 500ms gaps
# Row 40
$ws.Cells.Item(40, 6).ClearContents()
$ws.Cells.Item(40, 7).ClearContents()
$ws.Cells.Item(40, 8).Value = 0

# Row 43
$ws.Cells.Item(43, 6).ClearContents()
$ws.Cells.Item(43, 7).ClearContents()
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 10).Value = 107.3333333333339

# Row 50
$ws.Cells.Item(50, 6).ClearContents()
$ws.Cells.Item(50, 7).ClearContents()
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 11).Value = 557.3333333333339

# Row 55
$ws.Cells.Item(55, 6).ClearContents()
$ws.Cells.Item(55, 7).ClearContents()
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(55, 10).Value = 514.6666666666679

# Row 75
$ws.Cells.Item(75, 6).ClearContents()
$ws.Cells.Item(75, 7).ClearContents()
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 11).Value = 0
